$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new row of trade data as row 10, mirroring the structure of the
# preceding rows (A:H).
$ws.Cells.Item(10, 1).Value = 9822.0499999999993
$ws.Cells.Item(10, 2).Value = 9895.27
$ws.Cells.Item(10, 3).Value = 286
$ws.Cells.Item(10, 4).Value = 283.87
$ws.Cells.Item(10, 5).Value = $false
$ws.Cells.Item(10, 6).Value = -0.74
$ws.Cells.Item(10, 7).Value = 42612.673032407409
$ws.Cells.Item(10, 8).Value = $false

# Column G uses a date/time display format throughout the sheet; copy the
# formatting from the cell directly above so the new row matches exactly
# (reuses the existing date-time cell style rather than creating a new one).
$ws.Range("G9").Copy()
$ws.Range("G10").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false
